$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Parallel Script" runtime test was rewritten to cover different
# numbers of processors, so relabel the column header (C1) that feeds the
# chart series to make clear this run used 4 processors.
$ws.Range("C1").Value = "Parallel Script (4 Processors)"

# Move the active selection to D1 (reflecting where the user left off).
$ws.Range("D1").Select()
